$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.305.91"
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = "'1.855.86"
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'239.50"
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').Value = "'0.6960"
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = "'0.3079"
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = "'0.07560"
$ws.Range('E9').Value = '  +4.46%  '
$ws.Range('E10').Value = '  -3.82%  '
$ws.Range('D11').Value = "'0.08124"
$ws.Range('E11').Value = '  -3.13%  '
$ws.Range('D12').Value = "'1.874.41"
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = "'0.7267"
$ws.Range('E13').Value = '  -3.46%  '
$ws.Range('D14').Value = "'5.211"
$ws.Range('E14').Value = '  -3.96%  '
$ws.Range('D15').Value = "'89.44"
$ws.Range('E15').Value = '  -3.34%  '
$ws.Range('D16').Value = "'29.449.78"
$ws.Range('D17').Value = "'5.905"
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('D18').Value = "'242.65"
$ws.Range('E18').Value = '  -3.25%  '
$ws.Range('D19').Value = "'0.000007770"
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('D20').Value = "'13.15"
$ws.Range('E20').Value = '  -3.29%  '
$ws.Range('D21').Value = "'1.002"
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = "'2.126.10"
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = "'7.627"
$ws.Range('E24').Value = '  -5.04%  '
$ws.Range('D25').Value = "'9.066"
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'162.47"
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.1469"
$ws.Range('E27').Value = '  -6.15%  '
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E29').Value = '  -4.69%  '
$ws.Range('D30').Value = "'1.406"
$ws.Range('E30').Value = '  -7.27%  '
$ws.Range('D31').Value = "'1.517"
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = "'4.426"
$ws.Range('E32').Value = '  -3.72%  '
$ws.Range('D33').Value = "'4.052"
$ws.Range('E33').Value = '  -5.43%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('E35').Value = '  -3.17%  '
$ws.Range('D36').Value = "'0.7198"
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('D37').Value = "'0.9979"
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = "'2.667"
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = "'0.01870"
$ws.Range('E39').Value = '  -4.88%  '
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').Value = "'0.8841"
$ws.Range('E41').Value = '  +3.30%  '
$ws.Range('D42').Value = "'0.4314"
$ws.Range('E42').Value = '  -5.02%  '
$ws.Range('D43').Value = "'5.892"
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('D45').Value = "'1.047.47"
$ws.Range('E45').Value = '  -6.02%  '
$ws.Range('D46').Value = "'1.002"
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = "'102.69"
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = "'7.286"
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = "'2.017.92"
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'1.751"
$ws.Range('E50').Value = '  -5.67%  '
$ws.Range('D51').Value = "'9.271"
$ws.Range('E51').Value = '  -2.37%  '
